# Grace 4th updates on 0804
#
# 1) Two tiny floating point precision corrections on sheet "w0" (sheet5.xml)
# 2) Populate previously-empty sheet "s0" (sheet6.xml) with a male/female
#    age-profile table (rows for ages 6..22)

$wb = $excel.ActiveWorkbook

# --- 1) Precision fixes on sheet "w0" ----------------------------------
$wsW = $wb.Worksheets.Item("w0")
$wsW.Range("C10").Value = 0.0426580069922999
$wsW.Range("C25").Value = 0.0409819192244999

# --- 2) Fill sheet "s0" with data ---------------------------------------
$wsS = $wb.Worksheets.Item("s0")

$wsS.Range("B1").Value = "'male"
$wsS.Range("C1").Value = "'female"

$wsS.Range("A2").Value = "'6"
$wsS.Range("B2").Value = 1.35450501543214
$wsS.Range("C2").Value = 1.2854101799124

$wsS.Range("A3").Value = "'7"
$wsS.Range("B3").Value = 1.3471360470542
$wsS.Range("C3").Value = 1.27720069906436

$wsS.Range("A4").Value = "'8"
$wsS.Range("B4").Value = 1.28910228896852
$wsS.Range("C4").Value = 1.22257034104254

$wsS.Range("A5").Value = "'9"
$wsS.Range("B5").Value = 1.23106853088284
$wsS.Range("C5").Value = 1.16793998302072

$wsS.Range("A6").Value = "'10"
$wsS.Range("B6").Value = 1.17303477279716
$wsS.Range("C6").Value = 1.1133096249989

$wsS.Range("A7").Value = "'11"
$wsS.Range("B7").Value = 1.11500101471148
$wsS.Range("C7").Value = 1.05867926697708

$wsS.Range("A8").Value = "'12"
$wsS.Range("B8").Value = 1.0090941968288
$wsS.Range("C8").Value = 0.94635732204852

$wsS.Range("A9").Value = "'13"
$wsS.Range("B9").Value = 0.98844046300896
$wsS.Range("C9").Value = 0.92730304098396

$wsS.Range("A10").Value = "'14"
$wsS.Range("B10").Value = 0.96778672918912
$wsS.Range("C10").Value = 0.9082487599194

$wsS.Range("A11").Value = "'15"
$wsS.Range("B11").Value = 0.731669430135992
$wsS.Range("C11").Value = 0.67935702168538

$wsS.Range("A12").Value = "'16"
$wsS.Range("B12").Value = 0.715714220331216
$wsS.Range("C12").Value = 0.66479928438246

$wsS.Range("A13").Value = "'17"
$wsS.Range("B13").Value = 0.69975901052644
$wsS.Range("C13").Value = 0.65024154707954

$wsS.Range("A14").Value = "'18"
$wsS.Range("B14").Value = 0.338019695363148
$wsS.Range("C14").Value = 0.334515059293784

$wsS.Range("A15").Value = "'19"
$wsS.Range("B15").Value = 0.335517153359656
$wsS.Range("C15").Value = 0.332603630841872

$wsS.Range("A16").Value = "'20"
$wsS.Range("B16").Value = 0.333014611356163
$wsS.Range("C16").Value = 0.33069220238996

$wsS.Range("A17").Value = "'21"
$wsS.Range("B17").Value = 0.330512069352671
$wsS.Range("C17").Value = 0.328780773938047

$wsS.Range("A18").Value = "'22"
$wsS.Range("B18").Value = 0.328009527349178
$wsS.Range("C18").Value = 0.326869345486135
